# Auto-generated: applies scheduled-runner market data refresh to Pandaemonium_Profits sheets.
# Each sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) gets updated currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) for the rows whose source market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5551895.5
$ws.Range("I33").Value = 8327728.5
$ws.Range("J33").Value = 228.66667
$ws.Range("K33").Value = 8327728.5
$ws.Range("L33").Value = 228.66667
$ws.Range("M33").Value = -8327499.5
$ws.Range("N33").Value = -686.6666700000001
$ws.Range("H98").Value = 3591.3667
$ws.Range("I98").Value = 3521.64
$ws.Range("J98").Value = 3940
$ws.Range("K98").Value = 3521.64
$ws.Range("L98").Value = 3940
$ws.Range("M98").Value = -2023.64
$ws.Range("N98").Value = -6936
$ws.Range("H122").Value = 3591.3667
$ws.Range("I122").Value = 3521.64
$ws.Range("J122").Value = 3940
$ws.Range("K122").Value = 10564.92
$ws.Range("L122").Value = 11820
$ws.Range("M122").Value = -8114.92
$ws.Range("N122").Value = -16720
$ws.Range("H135").Value = 100001690
$ws.Range("I135").Value = 71429720
$ws.Range("J135").Value = 125002160
$ws.Range("K135").Value = 642867480
$ws.Range("L135").Value = 1125019440
$ws.Range("M135").Value = -642864945
$ws.Range("N135").Value = -1125024510
$ws.Range("H137").Value = 879392
$ws.Range("I137").Value = 2923.15
$ws.Range("J137").Value = 1853246.4
$ws.Range("K137").Value = 8769.450000000001
$ws.Range("L137").Value = 5559739.199999999
$ws.Range("M137").Value = -6219.450000000001
$ws.Range("N137").Value = -5564839.199999999
$ws.Range("H138").Value = 2711.3489
$ws.Range("I138").Value = 1419.2106
$ws.Range("J138").Value = 3734.2917
$ws.Range("K138").Value = 4257.6318
$ws.Range("L138").Value = 11202.8751
$ws.Range("M138").Value = 882.3681999999999
$ws.Range("N138").Value = -21482.8751
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3795.9524
$ws.Range("I132").Value = 3576.0625
$ws.Range("J132").Value = 4499.6
$ws.Range("K132").Value = 10728.1875
$ws.Range("L132").Value = 13498.8
$ws.Range("M132").Value = -8198.1875
$ws.Range("N132").Value = -18558.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 715788.7
$ws.Range("I31").Value = 5749.4287
$ws.Range("J31").Value = 1196783.1
$ws.Range("K31").Value = 5749.4287
$ws.Range("L31").Value = 1196783.1
$ws.Range("M31").Value = -5454.4287
$ws.Range("N31").Value = -1197373.1
$ws.Range("H34").Value = 715788.7
$ws.Range("I34").Value = 5749.4287
$ws.Range("J34").Value = 1196783.1
$ws.Range("K34").Value = 5749.4287
$ws.Range("L34").Value = 1196783.1
$ws.Range("M34").Value = -5547.4287
$ws.Range("N34").Value = -1197187.1
$ws.Range("H62").Value = 3279.8572
$ws.Range("I62").Value = 3326.5
$ws.Range("K62").Value = 3326.5
$ws.Range("M62").Value = -2702.5
$ws.Range("H65").Value = 3279.8572
$ws.Range("I65").Value = 3326.5
$ws.Range("K65").Value = 16632.5
$ws.Range("M65").Value = -13512.5
$ws.Range("H105").Value = 805.35297
$ws.Range("I105").Value = 694.2857
$ws.Range("K105").Value = 694.2857
$ws.Range("M105").Value = 1052.7143
$ws.Range("H132").Value = 4381.391
$ws.Range("I132").Value = 4265.222
$ws.Range("K132").Value = 12795.666
$ws.Range("M132").Value = -10265.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 389.73685
$ws.Range("J23").Value = 430.2353
$ws.Range("L23").Value = 1290.7059
$ws.Range("N23").Value = -1760.7059
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 9000
$ws.Range("N62").Value = -10372
$ws.Range("H63").Value = 2778.125
$ws.Range("I63").Value = 1665
$ws.Range("J63").Value = 4633.3335
$ws.Range("K63").Value = 4995
$ws.Range("L63").Value = 13900.0005
$ws.Range("M63").Value = -4246
$ws.Range("N63").Value = -15398.0005
$ws.Range("H64").Value = 1000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 27000
$ws.Range("N65").Value = -33864
$ws.Range("H66").Value = 2778.125
$ws.Range("I66").Value = 1665
$ws.Range("J66").Value = 4633.3335
$ws.Range("K66").Value = 14985
$ws.Range("L66").Value = 41700.0015
$ws.Range("M66").Value = -11241
$ws.Range("N66").Value = -49188.0015
$ws.Range("H67").Value = 1000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 2463.4243
$ws.Range("I68").Value = 896.34375
$ws.Range("J68").Value = 3938.3235
$ws.Range("K68").Value = 2689.03125
$ws.Range("L68").Value = 11814.9705
$ws.Range("M68").Value = -1878.03125
$ws.Range("N68").Value = -13436.9705
$ws.Range("H69").Value = 166667900
$ws.Range("I69").Value = 700
$ws.Range("J69").Value = 250001500
$ws.Range("K69").Value = 2100
$ws.Range("L69").Value = 750004500
$ws.Range("M69").Value = -1289
$ws.Range("N69").Value = -750006122
$ws.Range("H70").Value = 4666.6665
$ws.Range("J70").Value = 6000
$ws.Range("L70").Value = 18000
$ws.Range("N70").Value = -18630
$ws.Range("H71").Value = 2463.4243
$ws.Range("I71").Value = 896.34375
$ws.Range("J71").Value = 3938.3235
$ws.Range("K71").Value = 8067.09375
$ws.Range("L71").Value = 35444.9115
$ws.Range("M71").Value = -4011.09375
$ws.Range("N71").Value = -43556.9115
$ws.Range("H72").Value = 166667900
$ws.Range("I72").Value = 700
$ws.Range("J72").Value = 250001500
$ws.Range("K72").Value = 6300
$ws.Range("L72").Value = 2250013500
$ws.Range("M72").Value = -2244
$ws.Range("N72").Value = -2250021612
$ws.Range("H73").Value = 4666.6665
$ws.Range("J73").Value = 6000
$ws.Range("L73").Value = 18000
$ws.Range("N73").Value = -20184
$ws.Range("H88").Value = 3491.2
$ws.Range("J88").Value = 3491.2
$ws.Range("L88").Value = 10473.6
$ws.Range("N88").Value = -11329.6
$ws.Range("H91").Value = 3491.2
$ws.Range("J91").Value = 3491.2
$ws.Range("L91").Value = 10473.6
$ws.Range("N91").Value = -13437.6
$ws.Range("H107").Value = 592.7059
$ws.Range("J107").Value = 2542.8572
$ws.Range("L107").Value = 7628.571599999999
$ws.Range("N107").Value = -11468.5716
$ws.Range("H131").Value = 1213.4615
$ws.Range("I131").Value = 481
$ws.Range("J131").Value = 1466.0344
$ws.Range("K131").Value = 1443
$ws.Range("L131").Value = 4398.1032
$ws.Range("M131").Value = 3597
$ws.Range("N131").Value = -14478.1032

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3458.348
$ws.Range("I102").Value = 2879.5715
$ws.Range("J102").Value = 4358.6665
$ws.Range("K102").Value = 2879.5715
$ws.Range("L102").Value = 4358.6665
$ws.Range("M102").Value = -1257.5715
$ws.Range("N102").Value = -7602.6665
$ws.Range("H132").Value = 7740.5835
$ws.Range("I132").Value = 3749.7778
$ws.Range("J132").Value = 10135.066
$ws.Range("K132").Value = 11249.3334
$ws.Range("L132").Value = 30405.198
$ws.Range("M132").Value = -8719.3334
$ws.Range("N132").Value = -35465.198

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3549.0952
$ws.Range("I132").Value = 2271.4285
$ws.Range("J132").Value = 4187.9287
$ws.Range("K132").Value = 6814.2855
$ws.Range("L132").Value = 12563.7861
$ws.Range("M132").Value = -4284.2855
$ws.Range("N132").Value = -17623.7861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 52666.668
$ws.Range("J40").Value = 52666.668
$ws.Range("L40").Value = 52666.668
$ws.Range("N40").Value = -52964.668
$ws.Range("H46").Value = 39886
$ws.Range("J46").Value = 39886
$ws.Range("L46").Value = 39886
$ws.Range("N46").Value = -40348
$ws.Range("H126").Value = 1829.1765
$ws.Range("I126").Value = 1824
$ws.Range("J126").Value = 1841.6
$ws.Range("K126").Value = 5472
$ws.Range("L126").Value = 5524.799999999999
$ws.Range("M126").Value = -3002
$ws.Range("N126").Value = -10464.8
$ws.Range("H132").Value = 2981.1155
$ws.Range("I132").Value = 2755.5789
$ws.Range("J132").Value = 3593.2856
$ws.Range("K132").Value = 8266.736699999999
$ws.Range("L132").Value = 10779.8568
$ws.Range("M132").Value = -5736.736699999999
$ws.Range("N132").Value = -15839.8568
$ws.Range("H134").Value = 39886
$ws.Range("J134").Value = 39886
$ws.Range("L134").Value = 119658
$ws.Range("N134").Value = -124728
$ws.Range("H136").Value = 5350.4224
$ws.Range("I136").Value = 2208.2632
$ws.Range("J136").Value = 7646.615
$ws.Range("K136").Value = 6624.7896
$ws.Range("L136").Value = 22939.845
$ws.Range("M136").Value = -4074.7896
$ws.Range("N136").Value = -28039.845
$ws.Range("H138").Value = 38926.832
$ws.Range("J138").Value = 38926.832
$ws.Range("L138").Value = 38926.832
$ws.Range("N138").Value = -49206.832

